$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 <- values formerly belonging to row 10 (P8 / Lokalnamn stays unchanged)
$ws.Range("A8").Value = 111702486
$ws.Range("B8").Value = 90678
$ws.Range("E8").Value = 4366
$ws.Range("F8").Value = "Skarp dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum peckii"
$ws.Range("H8").Value = "Banker"
$ws.Range("Q8").Value = 517080.8398438052
$ws.Range("R8").Value = 6574959.907818918

# Row 10 <- values formerly belonging to row 11
$ws.Range("A10").Value = 111702506
$ws.Range("B10").Value = 90687
$ws.Range("E10").Value = 5964
$ws.Range("F10").Value = "Fjällig taggsvamp s.str."
$ws.Range("G10").Value = "Sarcodon imbricatus s.str."
$ws.Range("H10").Value = "(L.:Fr.) P.Karst."
$ws.Range("P10").Value = "Kyrkogården, Nrk"
$ws.Range("Q10").Value = 517093.6249861007
$ws.Range("R10").Value = 6574959.965416327

# Row 11 <- values formerly belonging to row 8
$ws.Range("A11").Value = 111702393
$ws.Range("B11").Value = 89183
$ws.Range("E11").Value = 3215
$ws.Range("F11").Value = "Rödgul trumpetsvamp"
$ws.Range("G11").Value = "Craterellus lutescens"
$ws.Range("H11").Value = "(Fr.) Fr."
$ws.Range("P11").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q11").Value = 517070.2129045375
$ws.Range("R11").Value = 6574934.844418272
